# Re-generate CPU summary table for supplementary_table_s5:
# The row for SBS_set2 / SigProfilerExtractor / seed.1076753 (row 77) is removed
# entirely (shifting all subsequent rows up by one), and the remaining
# SigProfilerExtractor cpu_time values are replaced with the newly measured,
# much shorter CPU times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 77 (SBS_set2 / SigProfilerExtractor / seed.1076753),
# shifting every row below it up by one.
$ws.Rows.Item(77).Delete()

# After the shift, rows 77-80 hold the remaining SigProfilerExtractor /
# SBS_set2 entries (seed.145879, seed.200437, seed.310111, seed.528401).
# Update their cpu_time (column D) values with the new measurements.
$ws.Range("D77").Value = 2102739.36
$ws.Range("D78").Value = 2107598.81
$ws.Range("D79").Value = 2086986.95
$ws.Range("D80").Value = 2113732.91
